# Mise à jour de l'application
#
# Adds a new training-session column (X) for the session dated 45874
# (2025-08-05), with one attendance marker per player, mirroring the
# formatting already used for the previous session column (W).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Attendance marker for each player row (row 1 is the header).
$values = @{
    2  = "P"   # Alban Rambaud
    3  = "P"   # Jassim Assoul
    4  = "P"   # Enzo Vita
    5  = "P"   # Romain Thunet
    6  = "P"   # Amine Taiar
    7  = "RH"  # Naim Ighbane
    8  = "P"   # Hedi Nasri
    9  = "P"   # Mattheo Haon
    10 = "P"   # Maé Clavel
    11 = "P"   # Levy Ndoutoume
    12 = "P"   # Yanis Berrached
    13 = "P"   # Rayane Chayebi
    14 = "P"   # Ilan Ihaddadene
    15 = "P"   # Karahali Souaré
    16 = "P"   # Amir Etien
    17 = "A"   # Karim Belmahi
    18 = "P"   # Emmanuel Valey
    19 = "P"   # Jeremie Laurent
    20 = "P"   # Sofiane Belle
    21 = "P"   # Amir Kherrab
    22 = "P"   # Naim Dhib
    23 = "P"   # Wael Fareh
    24 = "P"   # Yoan Zouma
    25 = "M"   # Ilyes Bougahnmi
    26 = "P"   # Omar Benyounes
    27 = "RH"  # Yoann Martelat
}

# Write the data cells first (with the same centered-text formatting used
# throughout column W) so every COUNTA/COUNTIF dependent on these rows
# recalculates immediately.
foreach ($row in $values.Keys) {
    $dst = "X" + $row
    $ws.Range($dst).Value = $values[$row]
    $ws.Range($dst).HorizontalAlignment = -4108
}

# Header cell (row 1): new session date, copying W1's date format/alignment.
$ws.Range("W1").Copy()
$ws.Range("X1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("X1").Value = 45874

$excel.Calculate()

# Restore (best-effort) the on-screen selection from the authored workbook.
$ws.Range("X1").Select()
$excel.ActiveWindow.ScrollColumn = 23
$ws.Range("AE13").Select()
